# Applies the CNV_metrics.docx table edit:
#  - inserts a new first column ("Metric" header + per-row metric names)
#  - resizes all 5 columns to 1728 dxa (from 4 columns at 2160 dxa)
#  - centers every cell's paragraph and sets its run font to
#    Times New Roman, 12pt (sz 24 half-points)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Insert a new first column (table goes from 4 -> 5 columns) ---
$firstCol = $t.Columns.Item(1)
$firstCol.Select()
$t.Columns.Add($firstCol)

# --- Resize every column to 1728 dxa (COM Width is expressed in points:
#     1728 twips / 20 = 86.4 pt) ---
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = 86.4
}

# --- Fill in the new first column's text for each row ---
$metricLabels = @(
    "Metric",
    "Percent difference between observed and expected coefficient of variation (2 Mbp window)",
    "Percent difference between observed and expected coefficient of variation (6 Mbp window)",
    "Correlation with label density",
    "Wave template correlation"
)

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Cell($i, 1)
    $cell.Range.Text = $metricLabels[$i - 1]
}

# --- Apply centered alignment + Times New Roman 12pt to every cell ---
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    for ($j = 1; $j -le $t.Columns.Count; $j++) {
        $cell = $t.Cell($i, $j)
        $p = $cell.Range.Paragraphs.Item(1)
        $p.Alignment = 1

        # Build a fresh Range (excluding the trailing cell-mark character)
        # instead of mutating cell.Range.End in place, which only retargets
        # formatting onto the tail of the text.
        $cellStart = $cell.Range.Start
        $cellEnd = $cell.Range.End
        $r = $d.Range($cellStart, $cellEnd - 1)
        $r.Font.Name = "Times New Roman"
        $r.Font.Size = 12
    }
}
